# qrCode now displays ticket owner
#
# The reservations log is updated:
#  - a new pending reservation for "monji" appears at row 4 (was "ahmed")
#  - "ahmed" is pushed down to row 5 and is now confirmed with an updated date
#  - the QR-scan confirmation log (previously rows 6-13 for "karoui") is
#    refreshed with new scan timestamps; the list is one entry shorter
#    (rows 6-12), so the sheet now spans A1:C12 instead of A1:C13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    # Writing date-like strings (e.g. "2025-02-26") straight into .Value lets
    # Excel auto-convert them into date serial numbers. Route them through a
    # text formula and then flatten it back down to a plain value so the
    # cell keeps its original style and ends up as a literal shared string.
    $cell.Formula = '="' + $text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

# Row 4: new pending reservation for "monji"
$ws.Cells.Item(4, 1).Value = "monji"
Set-TextValue $ws.Cells.Item(4, 2) "2025-02-26"
$ws.Cells.Item(4, 3).Value = "En attente"

# Row 5: "ahmed" reservation, now confirmed with an updated date
$ws.Cells.Item(5, 1).Value = "ahmed"
Set-TextValue $ws.Cells.Item(5, 2) "2025-02-25"
$ws.Cells.Item(5, 3).Value = "Confirmée"

# Rows 6-12: refreshed QR scan timestamps for "karoui" (status stays "Confirmé")
$ws.Cells.Item(6, 2).Value = "2025-03-06T13:56:23.644745200"
$ws.Cells.Item(7, 2).Value = "2025-03-06T13:58:54.494560500"
$ws.Cells.Item(8, 2).Value = "2025-03-06T14:04:51.406166900"
$ws.Cells.Item(9, 2).Value = "2025-03-06T14:07:54.549743800"
$ws.Cells.Item(10, 2).Value = "2025-03-06T14:10:38.794597700"
$ws.Cells.Item(11, 2).Value = "2025-03-06T14:19:21.469113600"
$ws.Cells.Item(12, 2).Value = "2025-03-06T14:24:27.672623800"

# Drop the oldest scan-log row; sheet now spans A1:C12 instead of A1:C13
$ws.Rows.Item(13).Delete()
